$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.049.13"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.259.09"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "2.600.86"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.863"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "2.251.25"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "43.923.50"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "0.0₃0990"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0855"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.115"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.81%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +23.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "1.804.36"
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.200"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "76.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
